# Add 2022-Q4 data:
#  - insert a new "2022-Q4" worksheet right after "总计", populated with the
#    newest quarter's fund-holding data (copied structure/styles from the
#    "2022-Q3" sheet, which keeps every other quarter sheet's own data/name
#    untouched - they just shift one tab to the right as a side effect)
#  - insert a new summary row into "总计" for the 2022-Q4 totals

$wb = $excel.ActiveWorkbook

$zj  = $wb.Worksheets.Item("总计")
$q3  = $wb.Worksheets.Item("2022-Q3")

# ---- 1. new "2022-Q4" sheet -------------------------------------------------
$q3.Copy([System.Reflection.Missing]::Value, $zj)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Columns B, D, E, G hold fund codes / figures that are stored as *text* in
# this workbook (e.g. "014038" keeps its leading zero). Force text format
# before writing so COM doesn't silently coerce them to numbers.
$q4.Range("B2:B3").NumberFormat = "@"
$q4.Range("D2:E3").NumberFormat = "@"
$q4.Range("G2:G3").NumberFormat = "@"

$q4.Cells.Item(2, 2).Value = "014038"
$q4.Cells.Item(2, 3).Value = "交银启诚混合A"
$q4.Cells.Item(2, 4).Value = "24.58"
$q4.Cells.Item(2, 5).Value = "81.04"
$q4.Cells.Item(2, 7).Value = "0.5752"
$q4.Cells.Item(2, 8).Value = 8

$q4.Cells.Item(3, 2).Value = "014039"
$q4.Cells.Item(3, 3).Value = "交银启诚混合C"
$q4.Cells.Item(3, 4).Value = "7.22"
$q4.Cells.Item(3, 5).Value = "81.04"
$q4.Cells.Item(3, 7).Value = "0.1689"
$q4.Cells.Item(3, 8).Value = 8

# ---- 2. new summary row in "总计" -------------------------------------------
# Extend the index-column (A) formatting down into row 6 first ...
$zj.Cells.Item(5, 1).Copy($zj.Cells.Item(6, 1))

# ... then push the existing quarter rows down by one (bottom-up so we never
# clobber a row before it has been copied).
for ($r = 5; $r -ge 2; $r--) {
    $zj.Cells.Item($r + 1, 2).Value = $zj.Cells.Item($r, 2).Value2
    $zj.Cells.Item($r + 1, 3).Value = $zj.Cells.Item($r, 3).Value2
    $zj.Cells.Item($r + 1, 4).Value = $zj.Cells.Item($r, 4).Value2
}
$zj.Cells.Item(6, 1).Value = 4

$zj.Cells.Item(2, 2).Value = "2022-Q4"
$zj.Cells.Item(2, 3).Value = 2
$zj.Cells.Item(2, 4).Value = 0.74
